$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Table 3 / Row 3 ("Поясніть, що іноді хтось..."):
#    - set the row's trHeight ("at least" rule)
#    - drop the 3 trailing empty paragraphs in cell 1
#    - drop the trailing empty paragraph + extra sentence in cell 2
# ----------------------------------------------------------------------
$row1 = $d.Tables.Item(3).Rows.Item(3)
$row1.HeightRule = 1
$row1.Height = 81.59033203125

$cell1a = $row1.Cells.Item(1)
for ($i = 0; $i -lt 3; $i++) {
    $n = $cell1a.Range.Paragraphs.Count
    $cell1a.Range.Paragraphs.Item($n).Range.Delete()
}

$cell1b = $row1.Cells.Item(2)
for ($i = 0; $i -lt 2; $i++) {
    $n = $cell1b.Range.Paragraphs.Count
    $cell1b.Range.Paragraphs.Item($n).Range.Delete()
}

# ----------------------------------------------------------------------
# 2) Table 3 / Row 10 ("Помічайте ознаки..." / "Вчіться самі"):
#    - translate "Помічайте ознаки, що щось не так " to English
#    - drop the 4 trailing paragraphs after "Зміни в настрої або поведінці"
# ----------------------------------------------------------------------
$d.Content.Find.Execute("Помічайте ознаки, що щось не так ", $true, $false, $false, $false, $false, $true, 1, $false, "Spot signs that something might be wrong ", 2) | Out-Null

$row2 = $d.Tables.Item(3).Rows.Item(10)
$cell2b = $row2.Cells.Item(2)
for ($i = 0; $i -lt 4; $i++) {
    $n = $cell2b.Range.Paragraphs.Count
    $cell2b.Range.Paragraphs.Item($n).Range.Delete()
}

# ----------------------------------------------------------------------
# 3) Translate "Покажіть, що дитина завжди може звернутися до вас.." line
# ----------------------------------------------------------------------
$d.Content.Find.Execute("Покажіть, що дитина завжди може звернутися до вас..", $true, $false, $false, $false, $false, $true, 1, $false, "Show them they can always come to you.", 2) | Out-Null

# ----------------------------------------------------------------------
# 4) Translate "Поясніть дитині, як приховувати чи видаляти публікації" line
# ----------------------------------------------------------------------
$d.Content.Find.Execute("Поясніть дитині, як приховувати чи видаляти публікації", $true, $false, $false, $false, $false, $true, 1, $false, "Hide or delete posts", 2) | Out-Null

Write-Output "done"
